# Updated cryptos list on Tue Apr 11 06:20:29 UTC 2023 with GitHub Actions
# Applies per-cell Coin/Link/Price/Volume(1h) updates to Sheet1 rows 2-51.
# Some Price values (column D) look numeric (e.g. "1.001", "0.4051"); they
# must stay as text like the original data, so those are entered with a
# leading apostrophe (forces text, same as typing '1.001 in Excel) and the
# cell style is reset back to Normal afterwards so no stray formatting
# (e.g. the quote-prefix indicator) is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.055.01'
$ws.Range('E2').Value = '  +5.44%  '
$ws.Range('D3').Value = '1.919.03'
$ws.Range('E3').Value = '  +2.33%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.78%  '
$ws.Range('D5').Value = "'327.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').Value = "'0.5252"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.30%  '
$ws.Range('D8').Value = "'0.4051"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.76%  '
$ws.Range('D9').Value = "'0.08470"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = "'1.128"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.22%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = "'42.86"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('D12').Value = "'22.22"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.81%  '
$ws.Range('D13').Value = "'6.356"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').Value = '1.918.95'
$ws.Range('E14').Value = '  +2.48%  '
$ws.Range('D15').Value = "'7.378"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = "'1.001"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'96.20"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.26%  '
$ws.Range('D18').Value = "'0.00001114"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').Value = "'0.06722"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = "'18.24"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.89%  '
$ws.Range('D21').Value = "'1.000"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = "'6.052"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').Value = '30.045.51'
$ws.Range('E23').Value = '  +5.33%  '
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('D25').Value = "'2.224"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').Value = '2.140.71'
$ws.Range('E26').Value = '  +2.49%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'21.14"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.52%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = "'160.07"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('D29').Value = "'2.456"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.74%  '
$ws.Range('D30').Value = "'129.41"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.79%  '
$ws.Range('D31').Value = "'1.082"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.06%  '
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').Value = "'6.096"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.55%  '
$ws.Range('D34').Value = "'3.660"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').Value = "'0.02520"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').Value = "'0.2232"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.30%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = "'5.216"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.23%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = "'1.236"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.66%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'9.006"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('D41').Value = "'0.6556"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = "'11.66"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.24%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'1.245"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = "'0.6190"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').Value = "'13.26"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('D46').Value = "'3.758"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').Value = "'2.072"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.30%  '
$ws.Range('D48').Value = "'125.81"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.32%  '
$ws.Range('E49').Value = '  +2.02%  '
$ws.Range('D50').Value = "'1.159"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('D51').Value = "'79.76"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.35%  '
